# Apply the "games updated 2019-05-14" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Résultats" sheet: enter game-3 scores for two series
#    (BOSTON vs CAROLINE, SAN JOSE vs ST-LOUIS)
# ---------------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Résultats")
$wsRes.Range("AC8").Value = 2
$wsRes.Range("AC9").Value = 1
$wsRes.Range("AB22").Value = 2
$wsRes.Range("AB23").Value = 4

# ---------------------------------------------------------------------------
# 2. "Pool" sheet: swap the two tied participants' prediction rows
#    (row 83 <-> row 87) so the sheet's pre-computed sort order settles.
# ---------------------------------------------------------------------------
$wsPool = $wb.Worksheets.Item("Pool")

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","T","U","V","W","X","Y","Z","AA","AC","AD","AE","AF")

foreach ($col in $cols) {
    $addr83 = $col + "83"
    $addr87 = $col + "87"
    $v83 = $wsPool.Range($addr83).Value2
    $v87 = $wsPool.Range($addr87).Value2
    $wsPool.Range($addr83).Value = $v87
    $wsPool.Range($addr87).Value = $v83
}

# ---------------------------------------------------------------------------
# 3. Restore the selections captured when the file was last saved.
# ---------------------------------------------------------------------------
$wsPool.Range("B6").Select()
$wsRes.Range("A2").Select()
